$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card text fields (rows 2-8) into a single python-tuple-like string
$combined = "('Lightning Dragon', ['{2}{R}{R}', 'Creature " + [char]0x2014 + " Dragon', 'Flying', 'Echo {2}{R}{R} (At the beginning of your upkeep, if this came under your control since the beginning of your last upkeep, sacrifice it unless you pay its echo cost.)', '{R}: Lightning Dragon gets +1/+0 until end of turn.', '4/4'])"

# Put the combined value in A2
$ws.Range("A2").Value = $combined

# Remove the now-unused rows 3-8
$ws.Range("A3:A8").EntireRow.Delete() | Out-Null
